# Biblioteca.xlsx - "USUARIOS" sheet: add library-user rows 2-9.
#
# The sheet started as a 1-row table (header + a single data row holding
# "2/2/2/2"). The edit grows it to 8 data rows:
#   - rows 2-8 -> all "1" placeholder values (A:D)
#   - row 9    -> the data that used to live in row 2, with the NOMBRE
#                 value bumped from "2" to "3"
# All values are plain numeric-looking strings that must stay text cells
# (mirrors the sheet's existing numberStoredAsText / ignoredError setup),
# so each write goes through a quick "format as text, assign, clear the
# format back off" round-trip instead of a bare .Value assignment (which
# Excel would otherwise silently coerce to a real number).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("USUARIOS")

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$row2 = @("1", "1", "1", "1")
for ($c = 1; $c -le 4; $c++) {
    Set-TextValue 2 $c $row2[$c - 1]
}

for ($r = 3; $r -le 8; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        Set-TextValue $r $c "1"
    }
}

$row9 = @("2", "3", "2", "2")
for ($c = 1; $c -le 4; $c++) {
    Set-TextValue 9 $c $row9[$c - 1]
}
